$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:ADSL"
$ws.Range("C2").Value = "NSE:GENSOL"
$ws.Range("D2").Value = "NSE:NYKAA"
$ws.Range("E2").Value = "NSE:CDSL"
$ws.Range("F2").Value = "NSE:NYKAA"

# Row 3
$ws.Range("B3").Value = "NSE:ARVINDFASN"
$ws.Range("C3").Value = "NSE:HEUBACHIND"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "NSE:HUDCO"
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("B4").Value = "NSE:BAJAJELEC"
$ws.Range("C4").Value = "NSE:HINDCOMPOS"
$ws.Range("E4").Value = "NSE:LT"

# Row 5
$ws.Range("B5").Value = "NSE:BATAINDIA"
$ws.Range("C5").Value = "NSE:LEXUS"
$ws.Range("E5").Value = "NSE:ONGC"

# Row 6
$ws.Range("B6").Value = "NSE:DEVYANI"
$ws.Range("C6").Value = "NSE:MADHUCON"

# Row 7
$ws.Range("B7").Value = "NSE:DPABHUSHAN"
$ws.Range("C7").Value = "NSE:MOKSH"

# Row 8
$ws.Range("B8").Value = "NSE:EICHERMOT"
$ws.Range("C8").Value = "NSE:NDL"

# Row 9
$ws.Range("B9").Value = "NSE:EIHAHOTELS"
$ws.Range("C9").Value = "NSE:PPL"

# Row 10
$ws.Range("B10").Value = "NSE:HATSUN"

# Row 11
$ws.Range("B11").Value = "NSE:INDIAGLYCO"

# Row 12
$ws.Range("B12").Value = "NSE:LEMONTREE"

# Row 13
$ws.Range("B13").Value = "NSE:NEOGEN"

# Row 14
$ws.Range("B14").Value = "NSE:NYKAA"

# Row 15
$ws.Range("B15").Value = "NSE:REDTAPE"

# Row 16
$ws.Range("B16").Value = "NSE:ROHLTD"

# Remove now-obsolete rows 17-19 (NSE:OSWALAGRO, NSE:PRESTIGE, NSE:RELAXO)
$ws.Rows("17:19").Delete()
